$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fv aksellast u 10t")

$ws.Range("B2").Value = 375523
$ws.Range("B3").Value = 379790
$ws.Range("B4").Value = 586656
$ws.Range("B5").Value = 792631
$ws.Range("B6").Value = 1643290
$ws.Range("B7").Value = 99643
$ws.Range("B8").Value = 532069
$ws.Range("B9").Value = 859737
$ws.Range("B10").Value = 1255210
$ws.Range("B11").Value = 1001920
